$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rewrite the "Reporte" field list (column E) and the "Diagnostico" field
#    list (column G). The Diagnostico list grows (gains tipo_accidente_id and
#    the brand-new "fallecio" field) and shifts down; the Reporte list loses
#    tipo_accidente_id and shifts up by one row.
# ---------------------------------------------------------------------------

# Clear the old column E / G content in the area that is being rearranged so
# no stale values are left behind once the tables are resized.
$ws.Range("E14:E25").ClearContents()
$ws.Range("G19:G32").ClearContents()

# New "Reporte" field list -> E13:E23
$reporte = @("Reporte","id_reporte","fecha","titular_minero","concesion","no_victimas","victima_nombre","empresa","tipo_empresa","no_decesos","observaciones")
for ($i = 0; $i -lt $reporte.Length; $i++) {
    $ws.Range("E" + (13 + $i)).Value2 = $reporte[$i]
}

# New "Diagnostico" field list -> G22:G37
$diagnostico = @("Diagnostico","id_diagnostico","sintomas","observaciones","oxigeno_sangre","pulsaciones","respiración","imagen","video","fecha","nombre_victima","ap_paterno_victima","ap_materno_victima","id_reporte","tipo_accidente_id","fallecio")
for ($i = 0; $i -lt $diagnostico.Length; $i++) {
    $ws.Range("G" + (22 + $i)).Value2 = $diagnostico[$i]
}

# ---------------------------------------------------------------------------
# 2. Resize the Excel Tables backing those two lists to their new extents.
# ---------------------------------------------------------------------------
$ws.ListObjects.Item("Tabla5").Resize($ws.Range("E13:E23"))
$ws.ListObjects.Item("Tabla9").Resize($ws.Range("G22:G37"))

# ---------------------------------------------------------------------------
# 3. Drawing clean-up: two of the connector arrows were exact duplicates of
#    others and got removed; the two remaining ones are stretched to reach
#    the relocated Diagnostico table.
# ---------------------------------------------------------------------------
$ws.Shapes.Item("Conector recto de flecha 4").Delete()
$ws.Shapes.Item("Conector recto de flecha 17").Delete()

$arrow1 = $ws.Shapes.Item("Conector recto de flecha 19")
$arrow1.Left = 194.17886841781495
$arrow1.Top = 190.50007874015748
$arrow1.Width = 340.5031502214567
$arrow1.Height = 341.16661417322837

$arrow2 = $ws.Shapes.Item("Conector recto de flecha 23")
$arrow2.Left = 453.0153257258858
$arrow2.Top = 204.16669291338582
$arrow2.Width = 82.50000000000006
$arrow2.Height = 315.8333858267717

# ---------------------------------------------------------------------------
# 4. View state: scroll/selection moved as the table grew past the bottom of
#    the window.
# ---------------------------------------------------------------------------
$ws.Range("H21").Select()
